$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9346511960029602
$ws.Range("B1").Value = 1.614339232444763
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.460383176803589
$ws.Range("E1").Value = 1.355513572692871
